# [RORO-5850][JE] Inspection UI: Display new Inspection location for DAERA Certex hold
#
# Replaces the "Manual Transit Procedure" (mtp) inspection heading translations
# (rows 2 & 3, KEY column + all language columns) with the new
# "DAERA CERTEX" (daera_certex) heading translations, and normalises the
# formatting of the data rows (font, row height) to the plain black
# Calibri/Arial look used by the rest of the sheet, dropping the old
# bespoke Helvetica Neue styling. The two residual, content-less bordered
# rows (4 & 5) underneath the table are also removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Row 2 — "inspection_needed_export.*.heading"
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "inspection_needed_export.daera_certex.heading"
$ws.Range("B2").Value = "For your DAERA CERTEX inspection"
$ws.Range("C2").Value = "Ar gyfer eich archwiliad DAERA CERTEX"
$ws.Range("D2").Value = "Do Twojej kontroli DAERA CERTEX"
$ws.Range("E2").Value = "Pentru inspecția dumneavoastră DAERA CERTEX"
$ws.Range("F2").Value = "Jūsų DAERA CERTEX patikrinimui"
$ws.Range("G2").Value = "За вашата CERTEX проверка от DAERA"
$ws.Range("H2").Value = "DAERA CERTEX-ellenőrzés esetén"
$ws.Range("I2").Value = "Para su inspección DAERA CERTEX"
$ws.Range("J2").Value = "Pour votre inspection DAERA CERTEX"
$ws.Range("K2").Value = "Für Ihre DAERA CERTEX-Prüfung"
$ws.Range("L2").Value = "K prohlídce DAERA CERTEX"
$ws.Range("M2").Value = "Za inspekciju CERTEX DAERA-e"

# ---------------------------------------------------------------------
# 2. Row 3 — "inspection_needed_import.*.heading" (same translated text
#    as row 2, only the KEY in column A differs)
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "inspection_needed_import.daera_certex.heading"
$ws.Range("B3").Value = "For your DAERA CERTEX inspection"
$ws.Range("C3").Value = "Ar gyfer eich archwiliad DAERA CERTEX"
$ws.Range("D3").Value = "Do Twojej kontroli DAERA CERTEX"
$ws.Range("E3").Value = "Pentru inspecția dumneavoastră DAERA CERTEX"
$ws.Range("F3").Value = "Jūsų DAERA CERTEX patikrinimui"
$ws.Range("G3").Value = "За вашата CERTEX проверка от DAERA"
$ws.Range("H3").Value = "DAERA CERTEX-ellenőrzés esetén"
$ws.Range("I3").Value = "Para su inspección DAERA CERTEX"
$ws.Range("J3").Value = "Pour votre inspection DAERA CERTEX"
$ws.Range("K3").Value = "Für Ihre DAERA CERTEX-Prüfung"
$ws.Range("L3").Value = "K prohlídce DAERA CERTEX"
$ws.Range("M3").Value = "Za inspekciju CERTEX DAERA-e"

# ---------------------------------------------------------------------
# 3. Re-style rows 2 & 3: plain black, non-wrapped, bottom-aligned text —
#    Arial for the KEY column (A), Calibri for every translation column.
# ---------------------------------------------------------------------
$dataRows = @(2, 3)
foreach ($r in $dataRows) {
    $ws.Rows.Item($r).RowHeight = 19.7

    $keyCell = $ws.Cells.Item($r, 1)
    $keyCell.Font.Name = "Arial"
    $keyCell.Font.Size = 16
    $keyCell.Font.Bold = $false
    $keyCell.Font.Color = 0x000000
    $keyCell.WrapText = $false
    $keyCell.VerticalAlignment = -4107

    for ($c = 2; $c -le 13; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Font.Name = "Calibri"
        $cell.Font.Size = 16
        $cell.Font.Bold = $false
        $cell.Font.Color = 0x000000
        $cell.WrapText = $false
        $cell.VerticalAlignment = -4107
    }
}

# ---------------------------------------------------------------------
# 4. Normalise the header row's theme-based black text to explicit black.
# ---------------------------------------------------------------------
$ws.Range("A1").Font.Color = 0x000000
$ws.Range("B1").Font.Color = 0x000000

# ---------------------------------------------------------------------
# 5. Drop the two empty, border-only rows below the table.
# ---------------------------------------------------------------------
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()
